$wb = $excel.ActiveWorkbook

# Sheet ALC row 118
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H118").Value = 2227
$ws.Range("I118").Value = 1341.3572
$ws.Range("J118").Value = 3260.25
$ws.Range("K118").Value = 4024.0716
$ws.Range("L118").Value = 9780.75
$ws.Range("M118").Value = -2367.0716
$ws.Range("N118").Value = -13094.75

# Sheet ALC row 125
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 1927.8572
$ws.Range("I125").Value = 2848.75
$ws.Range("J125").Value = 700
$ws.Range("K125").Value = 25638.75
$ws.Range("L125").Value = 6300
$ws.Range("M125").Value = -23178.75
$ws.Range("N125").Value = -11220

# Sheet ALC row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1074.84
$ws.Range("J129").Value = 1092.8125
$ws.Range("L129").Value = 3278.4375
$ws.Range("N129").Value = -13278.4375

# Sheet ALC row 134
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H134").Value = 55115
$ws.Range("J134").Value = 55115
$ws.Range("L134").Value = 55115
$ws.Range("N134").Value = -65255

# Sheet ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1797.9
$ws.Range("I2").Value = 1695.25
$ws.Range("J2").Value = 1866.3334
$ws.Range("K2").Value = 1695.25
$ws.Range("L2").Value = 1866.3334
$ws.Range("M2").Value = -1582.25
$ws.Range("N2").Value = -2092.3334

# Sheet ARM row 21
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 2500
$ws.Range("I21").Value = 2500
$ws.Range("K21").Value = 2500
$ws.Range("M21").Value = -2126

# Sheet ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7317
$ws.Range("I32").Value = 7435.174
$ws.Range("K32").Value = 7435.174
$ws.Range("M32").Value = -7148.174

# Sheet ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1797.9
$ws.Range("I116").Value = 1695.25
$ws.Range("J116").Value = 1866.3334
$ws.Range("K116").Value = 1695.25
$ws.Range("L116").Value = 1866.3334
$ws.Range("M116").Value = 598.75
$ws.Range("N116").Value = -6454.3334

# Sheet ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 13892645
$ws.Range("I132").Value = 22730728
$ws.Range("J132").Value = 4230.2856
$ws.Range("K132").Value = 68192184
$ws.Range("L132").Value = 12690.8568
$ws.Range("M132").Value = -68189654
$ws.Range("N132").Value = -17750.8568

# Sheet BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1797.9
$ws.Range("I3").Value = 1695.25
$ws.Range("J3").Value = 1866.3334
$ws.Range("K3").Value = 1695.25
$ws.Range("L3").Value = 1866.3334
$ws.Range("M3").Value = -1581.25
$ws.Range("N3").Value = -2094.3334

# Sheet BSM row 19
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()

# Sheet CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7250709.5
$ws.Range("I31").Value = 6515.7036
$ws.Range("J31").Value = 17545090
$ws.Range("K31").Value = 6515.7036
$ws.Range("L31").Value = 17545090
$ws.Range("M31").Value = -6220.7036
$ws.Range("N31").Value = -17545680

# Sheet CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 7250709.5
$ws.Range("I34").Value = 6515.7036
$ws.Range("J34").Value = 17545090
$ws.Range("K34").Value = 6515.7036
$ws.Range("L34").Value = 17545090
$ws.Range("M34").Value = -6313.7036
$ws.Range("N34").Value = -17545494

# Sheet CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2644.3044
$ws.Range("I58").Value = 1652.4286
$ws.Range("J58").Value = 4187.222
$ws.Range("K58").Value = 1652.4286
$ws.Range("L58").Value = 4187.222
$ws.Range("M58").Value = -1449.4286
$ws.Range("N58").Value = -4593.222

# Sheet CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4788.2144
$ws.Range("I132").Value = 4335.8887
$ws.Range("J132").Value = 5602.4
$ws.Range("K132").Value = 13007.6661
$ws.Range("L132").Value = 16807.2
$ws.Range("M132").Value = -10477.6661
$ws.Range("N132").Value = -21867.2

# Sheet CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2644.3044
$ws.Range("I136").Value = 1652.4286
$ws.Range("J136").Value = 4187.222
$ws.Range("K136").Value = 4957.2858
$ws.Range("L136").Value = 12561.666
$ws.Range("M136").Value = -2407.2858
$ws.Range("N136").Value = -17661.666

# Sheet CUL row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1201.6875
$ws.Range("I68").Value = 943.95
$ws.Range("J68").Value = 1385.7858
$ws.Range("K68").Value = 2831.85
$ws.Range("L68").Value = 4157.357400000001
$ws.Range("M68").Value = -2020.85
$ws.Range("N68").Value = -5779.357400000001

# Sheet CUL row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 1201.6875
$ws.Range("I71").Value = 943.95
$ws.Range("J71").Value = 1385.7858
$ws.Range("K71").Value = 8495.550000000001
$ws.Range("L71").Value = 12472.0722
$ws.Range("M71").Value = -4439.550000000001
$ws.Range("N71").Value = -20584.0722

# Sheet CUL row 118
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 823.0833
$ws.Range("I118").Value = 496.75
$ws.Range("J118").Value = 888.35
$ws.Range("K118").Value = 1490.25
$ws.Range("L118").Value = 2665.05
$ws.Range("M118").Value = -247.25
$ws.Range("N118").Value = -5151.05

# Sheet CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 740.59
$ws.Range("I131").Value = 408.21054
$ws.Range("J131").Value = 818.55554
$ws.Range("K131").Value = 1224.63162
$ws.Range("L131").Value = 2455.66662
$ws.Range("M131").Value = 3815.36838
$ws.Range("N131").Value = -12535.66662

# Sheet GSM row 10
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 2004
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 2004
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 2004
$ws.Range("M10").ClearContents()
$ws.Range("N10").Value = -2342

# Sheet GSM row 19
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 5329.75
$ws.Range("I19").Value = 319
$ws.Range("J19").Value = 7000
$ws.Range("K19").Value = 319
$ws.Range("L19").Value = 7000
$ws.Range("M19").Value = -31
$ws.Range("N19").Value = -7576

# Sheet GSM row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1000
$ws.Range("J113").Value = 1000
$ws.Range("L113").Value = 1000
$ws.Range("N113").Value = -5340

# Sheet GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5558330.5
$ws.Range("I122").Value = 11114278
$ws.Range("J122").Value = 2383.3333
$ws.Range("K122").Value = 33342834
$ws.Range("L122").Value = 7149.999899999999
$ws.Range("M122").Value = -33340384
$ws.Range("N122").Value = -12049.9999

# Sheet GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4082.468
$ws.Range("I132").Value = 3091.0645
$ws.Range("J132").Value = 6003.3125
$ws.Range("K132").Value = 9273.193499999999
$ws.Range("L132").Value = 18009.9375
$ws.Range("M132").Value = -6743.193499999999
$ws.Range("N132").Value = -23069.9375

# Sheet LTW row 3
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()

# Sheet LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6070
$ws.Range("I7").Value = 7309
$ws.Range("K7").Value = 7309
$ws.Range("M7").Value = -7197

# Sheet LTW row 15
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()

# Sheet LTW row 34
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value = 28449.75
$ws.Range("I34").Value = 27999.5
$ws.Range("J34").Value = 28900
$ws.Range("K34").Value = 27999.5
$ws.Range("L34").Value = 28900
$ws.Range("M34").Value = -27827.5
$ws.Range("N34").Value = -29244

# Sheet LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 13284.286
$ws.Range("I40").Value = 18663.334
$ws.Range("K40").Value = 18663.334
$ws.Range("M40").Value = -18527.334

# Sheet LTW row 43
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 25014
$ws.Range("J43").Value = 25014
$ws.Range("L43").Value = 25014
$ws.Range("N43").Value = -25400

# Sheet LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1477.7142
$ws.Range("I61").Value = 1547.8
$ws.Range("J61").Value = 1302.5
$ws.Range("K61").Value = 1547.8
$ws.Range("L61").Value = 1302.5
$ws.Range("M61").Value = -1345.8
$ws.Range("N61").Value = -1706.5

# Sheet LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1477.7142
$ws.Range("I113").Value = 1547.8
$ws.Range("J113").Value = 1302.5
$ws.Range("K113").Value = 1547.8
$ws.Range("L113").Value = 1302.5
$ws.Range("M113").Value = 622.2
$ws.Range("N113").Value = -5642.5

# Sheet LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 7504
$ws.Range("I122").Value = 14875
$ws.Range("J122").Value = 5236
$ws.Range("K122").Value = 44625
$ws.Range("L122").Value = 15708
$ws.Range("M122").Value = -42175
$ws.Range("N122").Value = -20608

# Sheet LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 6070
$ws.Range("I126").Value = 7309
$ws.Range("K126").Value = 21927
$ws.Range("M126").Value = -19457

# Sheet WVR row 5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 6666.3335
$ws.Range("J5").Value = 6666.3335
$ws.Range("L5").Value = 6666.3335
$ws.Range("N5").Value = -6890.3335

# Sheet WVR row 17
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 4
$ws.Range("I17").Value = 4
$ws.Range("K17").Value = 4
$ws.Range("M17").Value = 168

# Sheet WVR row 63
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 5226
$ws.Range("I63").Value = 5226
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 5226
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -4602
$ws.Range("N63").ClearContents()

# Sheet WVR row 66
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H66").Value = 5226
$ws.Range("I66").Value = 5226
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 15678
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -12558
$ws.Range("N66").ClearContents()
